$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4387992322444916
$ws.Range("B1").Value = 0.7151414155960083
$ws.Range("C1").Value = 2.179933547973633
$ws.Range("D1").Value = 4.72902774810791
$ws.Range("E1").Value = 2.181148290634155
